# Working on simulation - timeline issues.
# Apply the target edits to the workbook.

$wb = $excel.ActiveWorkbook

# --- Sheet "Model Parameters" : update ROUNDUP formula to use a literal
#     exponent of 10 instead of pulling the exponent from 'Indices Data'!B2.
$wsParams = $wb.Worksheets.Item("Model Parameters")
$wsParams.Range("B2").Formula = "=ROUNDUP(A2^10,0)"

# --- Sheet "Expected State Values" : refresh simulated output values.
$wsExp = $wb.Worksheets.Item("Expected State Values")

$wsExp.Range("C3").Value = 10.022
$wsExp.Range("F3").Value = 7.016
$wsExp.Range("I3").Value = 0.618
$wsExp.Range("M3").Value = 10.031

$wsExp.Range("F4").Value = 3.515
$wsExp.Range("R4").Value = 4.009

$wsExp.Range("F5").Value = 1.172
$wsExp.Range("R5").Value = 2.527

$wsExp.Range("F6").Value = 0.078
$wsExp.Range("R6").Value = 0.914

$wsExp.Range("F7").Value = 0.011
$wsExp.Range("R7").Value = 0.071

$wsExp.Range("R8").Value = 0.011

# Rows 13-22 (cols A-I) no longer hold placeholder data - clear them while
# keeping their existing styles.
$wsExp.Range("A13:I22").ClearContents()

$wsExp.Range("R15").Value = 2.528
$wsExp.Range("R16").Value = 0.914
$wsExp.Range("R17").Value = 0.242
$wsExp.Range("R18").Value = 0.07
$wsExp.Range("R19").Value = 0.011

$wsExp.Range("R26").Value = 0.914
$wsExp.Range("R27").Value = 0.241
$wsExp.Range("R28").Value = 0.07
$wsExp.Range("R29").Value = 0.011

$wsExp.Range("R37").Value = 0.241
$wsExp.Range("R38").Value = 0.07
$wsExp.Range("R39").Value = 0.011

$wsExp.Range("R48").Value = 0.07
$wsExp.Range("R49").Value = 0.011

$wsExp.Range("R59").Value = 0.011

# --- Selection / active-sheet bookkeeping.
# Move the "Expected State Values" selection to M9 (it is no longer the
# active tab once we switch to "Model Parameters" below).
$wsExp.Range("M9").Select()

# "Model Parameters" becomes the active (selected) sheet/tab.
$wsParams.Activate()
$wsParams.Range("B2").Select()

Write-Output "edits applied"
